$p = $ppt.ActivePresentation
$nl = [char]13

# -------------------------------------------------------------------------
# Slide 19 - "Understanding the PF Rules": merge two runs that made up one
# sentence into a single run (no visible text change, just de-fragmenting).
# -------------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$shRules = $s19.Shapes.Item(2)
$trRules = $shRules.TextFrame.TextRange
$paraDirection = $trRules.Paragraphs(2, 1)
$mergedSentence = "A matching pass rule in one direction automatically creates a matching pass rule in the other direction"
# Force an actual text-content change so the host collapses the paragraph
# back down to a single run (assigning the identical string is a no-op).
$paraDirection.Text = "."
$paraDirection.Text = $mergedSentence

# -------------------------------------------------------------------------
# Slide 30 - "Simplest way to organize the ruleset": rewrite the bullet
# list with new/expanded guidance.
# -------------------------------------------------------------------------
$s30 = $p.Slides.Item(30)
$shRuleset = $s30.Shapes.Item(2)

# Turn on "shrink text on overflow" for the body placeholder.
$shRuleset.TextFrame.AutoSize = 2

$trRuleset = $shRuleset.TextFrame.TextRange

# The trailing paragraph ("Allow the traffic and ports that you wish") is
# untouched by the edit, so anchor the rewrite to it and insert the new
# paragraphs immediately before it, then drop the three old paragraphs
# that used to precede it.
$anchorPara = $trRuleset.Paragraphs(4, 1)

$dash = [char]0x2013
$lq = [char]0x201C
$rq = [char]0x201D

$skipLine = "Skip the lo0 " + $dash + " set skip on lo0"
$blockLine = "Put the " + $lq + "block in all" + $rq + " at the top so that all traffic that you do not specify is dropped"
$passLine = "Allow traffic that your server initiates " + $lq + "pass out all" + $rq
$howeverLine = "However if the machine gets compromised, for example, is hacked and is used to attack other networks, you will have to change this rule"

$newParas = "Put your Macros at the top" + $nl + `
    "If you are on IPv6 add a macro for link local (fe80::/10)" + $nl + `
    $skipLine + $nl + `
    $blockLine + $nl + `
    $passLine + $nl + `
    $howeverLine + $nl + `
    "Allow link local addresses (if you are on an IPv6 network)" + $nl

$anchorPara.InsertBefore($newParas)

$trRuleset.Paragraphs(1, 1).Delete()
$trRuleset.Paragraphs(1, 1).Delete()
$trRuleset.Paragraphs(1, 1).Delete()

# Paragraph-level formatting: demote the two sub-bullets.
$trRuleset.Paragraphs(2, 1).IndentLevel = 2
$trRuleset.Paragraphs(6, 1).IndentLevel = 2

# --- Run-level splits -----------------------------------------------------

# P3: "Skip the lo0 - " / bold "set skip on lo0"
$p3 = $trRuleset.Paragraphs(3, 1)
$boldPrefixLen = ("Skip the lo0 " + $dash + " ").Length
$boldLen = "set skip on lo0".Length
$p3.Characters($p3.Start + $boldPrefixLen, $boldLen).Font.Bold = $true

# P4: "Put " / "the "block in all" at the top so that all traffic that you do not specify is dropped"
$p4 = $trRuleset.Paragraphs(4, 1)
$put4Len = "Put ".Length
$rest4Len = $p4.Length - $put4Len
$p4.Characters($p4.Start + $put4Len, $rest4Len).Text = $p4.Characters($p4.Start + $put4Len, $rest4Len).Text

# P5: "Allow traffic that your server initiates "pass " / "out" / " " / "all""
$p5 = $trRuleset.Paragraphs(5, 1)
$seg1 = "Allow traffic that your server initiates " + $lq + "pass "
$seg2 = "out"
$seg3 = " "
$seg4 = "all" + $rq
$off = $p5.Start
$p5.Characters($off + $seg1.Length, $seg2.Length).Text = $seg2
$p5.Characters($off + $seg1.Length + $seg2.Length, $seg3.Length).Text = $seg3
$p5.Characters($off + $seg1.Length + $seg2.Length + $seg3.Length, $seg4.Length).Text = $seg4

# P6: "However if ... change this " / "rule"
$p6 = $trRuleset.Paragraphs(6, 1)
$ruleLen = "rule".Length
$prefix6Len = $p6.Length - $ruleLen
$p6.Characters($p6.Start + $prefix6Len, $ruleLen).Text = "rule"

Write-Host "done"
